$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.622.82'
$ws.Range('E2').Value = '  +6.29%  '
$ws.Range('D3').Value = '2.294.34'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '305.04'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').Value = '101.24'
$ws.Range('E6').Value = '  +12.23%  '
$ws.Range('D7').Value = '0.567'
$ws.Range('E7').Value = '  +2.82%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.521'
$ws.Range('E9').Value = '  +6.36%  '
$ws.Range('D10').Value = '36.74'
$ws.Range('E10').Value = '  +11.60%  '
$ws.Range('E11').Value = '  +2.37%  '
$ws.Range('D12').Value = '7.37'
$ws.Range('E12').Value = '  +6.39%  '
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').Value = '2.645.15'
$ws.Range('E14').Value = '  +3.28%  '
$ws.Range('D15').Value = '2.294.27'
$ws.Range('E15').Value = '  +3.49%  '
$ws.Range('D16').Value = '13.90'
$ws.Range('E16').Value = '  +3.45%  '
$ws.Range('D17').Value = '0.815'
$ws.Range('E17').Value = '  +5.31%  '
$ws.Range('D18').Value = '46.616.41'
$ws.Range('D19').Value = '13.11'
$ws.Range('E19').Value = '  +13.49%  '
$ws.Range('D20').Value = '0.0₃0943'
$ws.Range('E20').Value = '  +4.82%  '
$ws.Range('D21').Value = '6.06'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('D22').Value = '66.50'
$ws.Range('E22').Value = '  +3.20%  '
$ws.Range('D23').Value = '248.14'
$ws.Range('E23').Value = '  +5.40%  '
$ws.Range('E24').Value = '  +4.19%  '
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('E26').Value = '  +3.85%  '
$ws.Range('D27').Value = '43.34'
$ws.Range('E27').Value = '  +13.57%  '
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  +2.02%  '
$ws.Range('E29').Value = '  +6.03%  '
$ws.Range('D30').Value = '20.04'
$ws.Range('E30').Value = '  +4.52%  '
$ws.Range('B31').Value = 'WEMIXToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D31').Value = '2.80'
$ws.Range('E31').Value = '  +11.82%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.68'
$ws.Range('E32').Value = '  +5.34%  '
$ws.Range('D33').Value = '147.49'
$ws.Range('E33').Value = '  -3.15%  '
$ws.Range('D34').Value = '0.0796'
$ws.Range('E34').Value = '  +5.55%  '
$ws.Range('D35').Value = '3.22'
$ws.Range('E35').Value = '  +13.22%  '
$ws.Range('E36').Value = '  +12.66%  '
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('E38').Value = '  +6.21%  '
$ws.Range('D39').Value = '16.07'
$ws.Range('E39').Value = '  +22.38%  '
$ws.Range('D40').Value = '4.05'
$ws.Range('E40').Value = '  +11.76%  '
$ws.Range('E41').Value = '  +8.06%  '
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').Value = '1.97'
$ws.Range('E44').Value = '  +11.78%  '
$ws.Range('D45').Value = '1.846.75'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('D46').Value = '87.28'
$ws.Range('E46').Value = '  +19.77%  '
$ws.Range('D47').Value = '0.196'
$ws.Range('E47').Value = '  +8.10%  '
$ws.Range('D48').Value = '73.92'
$ws.Range('E48').Value = '  +10.58%  '
$ws.Range('D49').Value = '4.91'
$ws.Range('E49').Value = '  +10.29%  '
$ws.Range('D50').Value = '95.96'
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').Value = '53.85'
$ws.Range('E51').Value = '  +6.16%  '
